$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text so values like
# "396.61" are not auto-converted to numbers, and reset the style
# back to Normal afterwards so no stray formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.796.69'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +10.49%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.251.56'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.82%  '

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '396.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.58%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.34'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.21%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.555'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.02%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.619'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.65%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.41%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0950'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +11.80%  '

$ws.Range("E12").Value = '  +2.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.792.34'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.67%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.99'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.55%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.04'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.55%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.263.29'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.04%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.04'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.57%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.75'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.82%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '56.748.72'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +10.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.29'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.93%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000105'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +9.60%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.83'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '303.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +14.60%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.72'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.83%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.18'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.12'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.42%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.93'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.06%  '

$ws.Range("E28").Value = '  +4.94%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.43%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.168'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.00%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.110'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.91%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.98'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '37.28'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0482'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.76%  '

$ws.Range("E36").Value = '  +3.07%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.48'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.78%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.14'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +24.36%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.54'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.45%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.14%  '

$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '133.96'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.44%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.92'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.23%  '

$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.31'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.83%  '

$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.119'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.81%  '

$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.96'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.05%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.279'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '21.91'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.77%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.150.09'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.70%  '

$ws.Range("E49").Value = '  +2.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.04'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +43.18%  '

$ws.Range("E51").Value = '  -3.83%  '
